$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "PIR" - append rows 68-80 (Bathroom / No Motion / Inactive)
# ---------------------------------------------------------------------------
$pir = $wb.Worksheets.Item("PIR")

$pirData = @(
    @("2026-01-30", "18:24:10", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:12", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:18", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:22", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:28", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:32", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:37", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:42", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:48", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:53", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:24:58", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:25:03", "18:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "18:25:08", "18:00", "Bathroom", "No Motion", "Inactive")
)

$startRow = 68
for ($i = 0; $i -lt $pirData.Count; $i++) {
    $row = $startRow + $i
    $values = $pirData[$i]

    # Column A holds a date-shaped string ("2026-01-30"); force text so Excel
    # doesn't coerce it into a date serial number.
    $cellA = $pir.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.ClearFormats()

    $pir.Cells.Item($row, 2).Value = $values[1]
    $pir.Cells.Item($row, 3).Value = $values[2]
    $pir.Cells.Item($row, 4).Value = $values[3]
    $pir.Cells.Item($row, 5).Value = $values[4]
    $pir.Cells.Item($row, 6).Value = $values[5]
}

# ---------------------------------------------------------------------------
# Sheet "Humidity" - append rows 49-55 (Bathroom / %RH / Active)
# ---------------------------------------------------------------------------
$humidity = $wb.Worksheets.Item("Humidity")

$humidityData = @(
    @("2026-01-30", "18:24:10", "18:00", "Bathroom", "86.7%", "Active"),
    @("2026-01-30", "18:24:13", "18:00", "Bathroom", "85.7%", "Active"),
    @("2026-01-30", "18:24:28", "18:00", "Bathroom", "86.7%", "Active"),
    @("2026-01-30", "18:24:33", "18:00", "Bathroom", "86.7%", "Active"),
    @("2026-01-30", "18:24:48", "18:00", "Bathroom", "86.8%", "Active"),
    @("2026-01-30", "18:24:53", "18:00", "Bathroom", "85.8%", "Active"),
    @("2026-01-30", "18:25:08", "18:00", "Bathroom", "86.8%", "Active")
)

$startRow = 49
for ($i = 0; $i -lt $humidityData.Count; $i++) {
    $row = $startRow + $i
    $values = $humidityData[$i]

    # Column A holds a date-shaped string; force text.
    $cellA = $humidity.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $values[0]
    $cellA.ClearFormats()

    $humidity.Cells.Item($row, 2).Value = $values[1]
    $humidity.Cells.Item($row, 3).Value = $values[2]
    $humidity.Cells.Item($row, 4).Value = $values[3]

    # Column E holds a percentage-shaped string ("86.7%"); force text so
    # Excel doesn't coerce it into a numeric percentage value.
    $cellE = $humidity.Cells.Item($row, 5)
    $cellE.NumberFormat = "@"
    $cellE.Value = $values[4]
    $cellE.ClearFormats()

    $humidity.Cells.Item($row, 6).Value = $values[5]
}
